# Update the multiplication answers in the table to match the new
# output generated at c8c62b6.

$d = $word.ActiveDocument

$replacements = @(
    @{ old = "15×33=495";  new = "72×35=2520" },
    @{ old = "57×14=798";  new = "91×85=7735" },
    @{ old = "99×14=1386"; new = "58×90=5220" },
    @{ old = "14×74=1036"; new = "48×27=1296" },
    @{ old = "12×62=744";  new = "79×50=3950" },
    @{ old = "18×81=1458"; new = "34×35=1190" },
    @{ old = "90×72=6480"; new = "60×42=2520" },
    @{ old = "46×81=3726"; new = "42×15=630" },
    @{ old = "11×59=649";  new = "47×47=2209" },
    @{ old = "67×26=1742"; new = "98×31=3038" },
    @{ old = "66×77=5082"; new = "60×56=3360" },
    @{ old = "19×42=798";  new = "31×55=1705" },
    @{ old = "25×76=1900"; new = "39×81=3159" },
    @{ old = "85×12=1020"; new = "24×12=288" },
    @{ old = "33×76=2508"; new = "62×52=3224" },
    @{ old = "58×61=3538"; new = "16×50=800" },
    @{ old = "68×45=3060"; new = "77×25=1925" },
    @{ old = "92×23=2116"; new = "90×90=8100" },
    @{ old = "82×17=1394"; new = "88×46=4048" },
    @{ old = "37×87=3219"; new = "15×74=1110" },
    @{ old = "80×22=1760"; new = "86×29=2494" },
    @{ old = "50×91=4550"; new = "32×49=1568" },
    @{ old = "27×99=2673"; new = "53×18=954" },
    @{ old = "30×47=1410"; new = "42×67=2814" },
    @{ old = "39×34=1326"; new = "85×64=5440" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
